# Regenerate sval data to filter save games:
# update the numeric stats on row 2 (B2:E2, G2) of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.445647641019636
$ws.Range("C2").Value = 0.3048912486333797
$ws.Range("D2").Value = 0.7210945179870265
$ws.Range("E2").Value = 0.5333859586016987
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.005019366241741
